$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new entry "103. Binary Tree Zigzag Level Order Traversal" ---

# A10: Question No = 103 (keeps existing left/top/wrap style)
$ws.Range("A10").Value = 103

# B10: GFG/LC = "LC" (keeps existing left/top/wrap style)
$ws.Range("B10").Value = "LC"

# C10: Question text - reset formatting back to the default "Normal" style
# (this cell previously inherited the row's left/top/wrap style, the new
# row no longer carries that formatting on column C, matching row 9's C9)
$ws.Range("C10").Value = "Binary Tree Zigzag Level Order Traversal"
$ws.Range("C10").Style = "Normal"

# D10: Java/Python - vertical-top alignment only (same as D7/D8/D9)
$ws.Range("D10").Value = "Java/Python"
$ws.Range("D10").VerticalAlignment = -4160

# E10: Difficulty = "Medium", highlighted with a yellow fill
$ws.Range("E10").Value = "Medium"
$ws.Range("E10").Interior.Color = 65535

# --- View: scroll/selection update ---
$window = $excel.ActiveWindow
$window.ScrollRow = 7
$window.ScrollColumn = 2
$ws.Range("C16").Select()
